# CIERRE 11 JUN 22
# Update the "ARQUITECTO" vale sheet (the active/selected sheet): the paid
# amount drops from 150,000 to 50,000, and the amount-in-words cell is
# updated to match ("CIENTO CINCUENTA MIL..." -> "CIENTO MIL...").
# The user's selection also moved from C14 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1: vale amount 150000 -> 50000
$ws.Range("D1").Value = 50000

# A2: amount spelled out in words, updated to match the new amount
$ws.Range("A2").Value = "CIENTO    MIL   PESOS 00/100 M.N."

# Cursor/selection ends up on C3
$ws.Range("C3").Select()
